# Applies the "ModelCanvas" Raster LatLong calculation table to Sheet2,
# and nudges Sheet1's active-cell selection from D9 to E9.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: move the active selection ---
[void]$ws1.Range("E9").Select()

# --- Sheet2: build the Raster LatLongBox table ---

# Column widths (match bestFit widths used in the authored workbook, as
# closely as the host's pixel-quantized column-width model allows)
$ws2.Columns.Item(1).ColumnWidth = 13.71
$ws2.Columns.Item(2).ColumnWidth = 21.71
$ws2.Columns.Item(3).ColumnWidth = 19.86
$ws2.Columns.Item(4).ColumnWidth = 22.71
$ws2.Columns.Item(5).ColumnWidth = 17.29
$ws2.Columns.Item(6).ColumnWidth = 14.71

# Header row
$ws2.Range("B1").Value = "Model"
$ws2.Range("C1").Value = "Raster:"
$ws2.Range("D1").Value = "Raster LatLongBox"
$ws2.Range("E1").Value = "Buffer LatLongBox"

# Row labels
$ws2.Range("A2").Value = "North"
$ws2.Range("A3").Value = "South"
$ws2.Range("A4").Value = "East"
$ws2.Range("A5").Value = "West"

# Data values
$ws2.Range("B2").Value = 44.323888888653777
$ws2.Range("C2").Value = 44.285831541747001
$ws2.Range("D2").Value = 44.285831541747001
$ws2.Range("E2").Value = 44.231388332300398

$ws2.Range("B3").Value = 44.211666666424001
$ws2.Range("C3").Value = 44.239627838134702
$ws2.Range("D3").Value = 44.239627838134702
$ws2.Range("E3").Value = 44.211758699770101

$ws2.Range("B4").Value = -71.244722218675577
$ws2.Range("C4").Value = -71.279353660393099
$ws2.Range("D4").Value = -71.279353660393099
$ws2.Range("E4").Value = -71.292219848815805

$ws2.Range("B5").Value = -71.384722218684999
$ws2.Range("C5").Value = -71.325927734375
$ws2.Range("D5").Value = -71.325927734375
$ws2.Range("E5").Value = -71.384719848632798

# Difference rows
$ws2.Range("B7").Formula = "=B2-B3"
$ws2.Range("C7").Formula = "=C2-C3"

$ws2.Range("B8").Formula = "=B4-B5"
$ws2.Range("C8").Formula = "=C4-C5"

# Tile size rows
$ws2.Range("A10").Value = "Tile Lat Height"
$ws2.Range("C10").Value = 0.092499999816936906

$ws2.Range("A11").Value = "Tile Lon Width"
$ws2.Range("C11").Value = 0.092499999816936906

# Number formats: B column uses 17 decimal places, C/D/E use 12 decimal places
$ws2.Range("B2:B5").NumberFormat = "0.00000000000000000"
$ws2.Range("B7:B8").NumberFormat = "0.00000000000000000"

$ws2.Range("C2:E5").NumberFormat = "0.000000000000"
$ws2.Range("C7:C8").NumberFormat = "0.000000000000"
$ws2.Range("C10").NumberFormat = "0.000000000000"
$ws2.Range("C11").NumberFormat = "0.000000000000"

# Active selection on Sheet2
[void]$ws2.Range("B4").Select()

# Leave Sheet1 as the active sheet/tab (tabSelected) as in the source file
[void]$ws1.Activate()
